# Applies: feat: enhance schedule loading with debug logging and improve CSS for special cells
#
# Concretely (from the OOXML diff):
#  - regseason sheet: deselect tab, change selection to header row (A1:XFD1),
#    and tidy a couple of column widths (split merged-width col pairs).
#  - playoffs sheet: populate with the same header row used on regseason and
#    copy the trailing 3 "playoff week" rows from regseason (weeks 32-34,
#    type P, with their formatted date cells); select the whole sheet.
#  - new "champs" sheet appended at the end (and made the active tab) holding
#    a year-by-year championship history table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) regseason: selection / tab state
# ---------------------------------------------------------------------
$regseason = $wb.Worksheets.Item("regseason")
$regseason.Range("A1:XFD1").Select() | Out-Null

# Column width tidy-up (matches the target cols metadata as closely as the
# host's AutoFit heuristic allows).
$regseason.Columns.Item(2).AutoFit() | Out-Null
$regseason.Columns.Item(8).AutoFit() | Out-Null
$regseason.Columns.Item(9).AutoFit() | Out-Null
$regseason.Columns.Item(14).AutoFit() | Out-Null
$regseason.Columns.Item(15).AutoFit() | Out-Null
$regseason.Columns.Item(20).AutoFit() | Out-Null
$regseason.Columns.Item(21).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 2) playoffs: header row + the 3 playoff-week rows copied from regseason
# ---------------------------------------------------------------------
$playoffs = $wb.Worksheets.Item("playoffs")

# Header row (identical to regseason's).
$regseason.Range("A1:U1").Copy() | Out-Null
$playoffs.Range("A1").PasteSpecial(-4104) | Out-Null

# Data: regseason rows 33-35 (weeks 32-34, "P" type w/ formatted dates)
# become playoffs rows 2-4.
$regseason.Range("A33:C35").Copy() | Out-Null
$playoffs.Range("A2").PasteSpecial(-4104) | Out-Null

# Re-stamp the date-formatted style (numFmtId 14) on the copied date cells -
# a plain "paste all" loses the number format style id, a formats-only paste
# restores it without disturbing the text values just pasted.
$regseason.Range("C33:C35").Copy() | Out-Null
$playoffs.Range("C2:C4").PasteSpecial(-4122) | Out-Null

$playoffs.Columns.Item(2).AutoFit() | Out-Null
$playoffs.Columns.Item(8).AutoFit() | Out-Null
$playoffs.Columns.Item(9).AutoFit() | Out-Null
$playoffs.Columns.Item(14).AutoFit() | Out-Null
$playoffs.Columns.Item(15).AutoFit() | Out-Null
$playoffs.Columns.Item(20).AutoFit() | Out-Null
$playoffs.Columns.Item(21).AutoFit() | Out-Null

$playoffs.Cells.Select() | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) champs: brand-new sheet, appended after playoffs, becomes active tab
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$champs = $wb.Sheets.Add($null, $lastSheet)
$champs.Name = "champs"

# Header row
$champs.Cells.Item(1,1).Value = "year"
$champs.Cells.Item(1,2).Value = "team"
$champs.Cells.Item(1,3).Value = "venue"
$champs.Cells.Item(1,4).Value = "p1"
$champs.Cells.Item(1,5).Value = "p2"
$champs.Cells.Item(1,6).Value = "p3"
$champs.Cells.Item(1,7).Value = "p4"
$champs.Cells.Item(1,8).Value = "p5"
$champs.Cells.Item(1,9).Value = "p6"
$champs.Cells.Item(1,10).Value = "p7"
$champs.Cells.Item(1,11).Value = "p8"
$champs.Cells.Item(1,12).Value = "p9"

# Row 2 - 2025 champions with full roster
$champs.Cells.Item(2,1).Value = 2025
$champs.Cells.Item(2,2).Value = "Team 2"
$champs.Cells.Item(2,3).Value = "Magna Centre"
$champs.Cells.Item(2,4).Value = "Frank Vucko"
$champs.Cells.Item(2,5).Value = "James Wang"
$champs.Cells.Item(2,6).Value = "Mav Marick"
$champs.Cells.Item(2,7).Value = "George Sparangis"
$champs.Cells.Item(2,8).Value = "Petar Rafajlovic"
$champs.Cells.Item(2,9).Value = "Ralph Romano"
$champs.Cells.Item(2,10).Value = "Richard Montoya"
$champs.Cells.Item(2,11).Value = "Sean Ludwig"

# Rows 3-10: Magna Centre years, team only
$magnaYears = @(
    @(2024, "Team 3"),
    @(2023, "Team 4"),
    @(2022, "Team 5"),
    @(2021, "Team 6"),
    @(2020, "Team 2"),
    @(2019, "Team 3"),
    @(2018, "Team 4"),
    @(2017, "Team 5")
)
$r = 3
foreach ($entry in $magnaYears) {
    $champs.Cells.Item($r,1).Value = $entry[0]
    $champs.Cells.Item($r,2).Value = $entry[1]
    $champs.Cells.Item($r,3).Value = "Magna Centre"
    $r = $r + 1
}

# Rows 11-21: Sacred Heart years, team only
$sacredYears = @(
    @(2016, "Team 6"),
    @(2015, "Team 2"),
    @(2014, "Team 3"),
    @(2013, "Team 4"),
    @(2012, "Team 5"),
    @(2011, "Team 6"),
    @(2010, "Team 2"),
    @(2009, "Team 3"),
    @(2008, "Team 4"),
    @(2007, "Team 5"),
    @(2006, "Team 6")
)
foreach ($entry in $sacredYears) {
    $champs.Cells.Item($r,1).Value = $entry[0]
    $champs.Cells.Item($r,2).Value = $entry[1]
    $champs.Cells.Item($r,3).Value = "Sacred Heart"
    $r = $r + 1
}

$champs.Columns.Item("A:L").AutoFit() | Out-Null
